$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-17 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-18 Friday", 2)

# Update the division problems in the table, addressed by (row, column)
# to avoid ambiguity from duplicate/cross-moved values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "94÷3=31, 1"
$t.Cell(1, 2).Range.Text = "48÷9=5, 3"
$t.Cell(1, 3).Range.Text = "62÷9=6, 8"
$t.Cell(1, 4).Range.Text = "68÷7=9, 5"
$t.Cell(1, 5).Range.Text = "18÷9=2, 0"
$t.Cell(5, 1).Range.Text = "21÷5=4, 1"
$t.Cell(5, 2).Range.Text = "84÷7=12, 0"
$t.Cell(5, 3).Range.Text = "68÷3=22, 2"
$t.Cell(5, 4).Range.Text = "66÷3=22, 0"
$t.Cell(5, 5).Range.Text = "36÷4=9, 0"
$t.Cell(9, 1).Range.Text = "34÷5=6, 4"
$t.Cell(9, 2).Range.Text = "65÷4=16, 1"
$t.Cell(9, 3).Range.Text = "21÷4=5, 1"
$t.Cell(9, 4).Range.Text = "12÷8=1, 4"
$t.Cell(9, 5).Range.Text = "48÷3=16, 0"
$t.Cell(13, 1).Range.Text = "70÷4=17, 2"
$t.Cell(13, 2).Range.Text = "89÷9=9, 8"
$t.Cell(13, 3).Range.Text = "13÷6=2, 1"
$t.Cell(13, 4).Range.Text = "33÷6=5, 3"
$t.Cell(13, 5).Range.Text = "54÷6=9, 0"
$t.Cell(17, 1).Range.Text = "50÷6=8, 2"
$t.Cell(17, 2).Range.Text = "33÷2=16, 1"
$t.Cell(17, 3).Range.Text = "70÷6=11, 4"
$t.Cell(17, 4).Range.Text = "84÷5=16, 4"
$t.Cell(17, 5).Range.Text = "25÷3=8, 1"
